$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
# Row 32
$ws.Range("H32").Value = 1175.25
$ws.Range("J32").Value = 1383.6666
$ws.Range("L32").Value = 1383.6666
$ws.Range("N32").Value = -2035.6666
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null
# Row 106
$ws.Range("H106").Value = 2446.3
$ws.Range("I106").Value = 2495.889
$ws.Range("K106").Value = 2495.889
$ws.Range("M106").Value = -1864.889
# Row 107
$ws.Range("H107").Value = 678.6923
$ws.Range("I107").Value = 567.4286
$ws.Range("J107").Value = 808.5
$ws.Range("K107").Value = 567.4286
$ws.Range("L107").Value = 808.5
$ws.Range("M107").Value = 1352.5714
$ws.Range("N107").Value = -4648.5
# Row 129
$ws.Range("H129").Value = 1181.5883
$ws.Range("J129").Value = 1277.6957
$ws.Range("L129").Value = 3833.0871
$ws.Range("N129").Value = -13833.0871
# Row 137
$ws.Range("H137").Value = 2771.56
$ws.Range("I137").Value = 2120
$ws.Range("J137").Value = 3600.818
$ws.Range("K137").Value = 6360
$ws.Range("L137").Value = 10802.454
$ws.Range("M137").Value = -3810
$ws.Range("N137").Value = -15902.454

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 20614.889
$ws.Range("I32").Value = 24342.678
$ws.Range("J32").Value = 7567.625
$ws.Range("K32").Value = 24342.678
$ws.Range("L32").Value = 7567.625
$ws.Range("M32").Value = -24055.678
$ws.Range("N32").Value = -8141.625
# Row 61
$ws.Range("H61").Value = 2300.6943
$ws.Range("I61").Value = 2176.5454
$ws.Range("J61").Value = 3666.3333
$ws.Range("K61").Value = 2176.5454
$ws.Range("L61").Value = 3666.3333
$ws.Range("M61").Value = -1964.5454
$ws.Range("N61").Value = -4090.3333
# Row 74
$ws.Range("H74").Value = 1240.431
$ws.Range("I74").Value = 1115.85
$ws.Range("J74").Value = 1517.2778
$ws.Range("K74").Value = 1115.85
$ws.Range("L74").Value = 1517.2778
$ws.Range("M74").Value = -241.8499999999999
$ws.Range("N74").Value = -3265.2778
# Row 77
$ws.Range("H77").Value = 1240.431
$ws.Range("I77").Value = 1115.85
$ws.Range("J77").Value = 1517.2778
$ws.Range("K77").Value = 5579.25
$ws.Range("L77").Value = 7586.389
$ws.Range("M77").Value = -1211.25
$ws.Range("N77").Value = -16322.389
# Row 122
$ws.Range("H122").Value = 4463.091
$ws.Range("I122").Value = 5015.4614
$ws.Range("J122").Value = 2411.4285
$ws.Range("K122").Value = 15046.3842
$ws.Range("L122").Value = 7234.2855
$ws.Range("M122").Value = -12596.3842
$ws.Range("N122").Value = -12134.2855
# Row 132
$ws.Range("H132").Value = 4315.478
$ws.Range("I132").Value = 4151.4883
$ws.Range("K132").Value = 12454.4649
$ws.Range("M132").Value = -9924.464899999999
# Row 136
$ws.Range("H136").Value = 2300.6943
$ws.Range("I136").Value = 2176.5454
$ws.Range("J136").Value = 3666.3333
$ws.Range("K136").Value = 6529.6362
$ws.Range("L136").Value = 10998.9999
$ws.Range("M136").Value = -3979.6362
$ws.Range("N136").Value = -16098.9999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 29261.666
$ws.Range("I7").Value = 28883
$ws.Range("J7").Value = 29451
$ws.Range("K7").Value = 28883
$ws.Range("L7").Value = 29451
$ws.Range("M7").Value = -28770
$ws.Range("N7").Value = -29677
# Row 86
$ws.Range("H86").Value = 79023.92
$ws.Range("I86").Value = 2098.75
$ws.Range("J86").Value = 202104.2
$ws.Range("K86").Value = 2098.75
$ws.Range("L86").Value = 202104.2
$ws.Range("M86").Value = -975.75
$ws.Range("N86").Value = -204350.2
# Row 89
$ws.Range("H89").Value = 79023.92
$ws.Range("I89").Value = 2098.75
$ws.Range("J89").Value = 202104.2
$ws.Range("K89").Value = 10493.75
$ws.Range("L89").Value = 1010521
$ws.Range("M89").Value = -4877.75
$ws.Range("N89").Value = -1021753
# Row 94
$ws.Range("H94").Value = 1122.0869
$ws.Range("I94").Value = 994.64703
$ws.Range("J94").Value = 1483.1666
$ws.Range("K94").Value = 994.64703
$ws.Range("L94").Value = 1483.1666
$ws.Range("M94").Value = -543.64703
$ws.Range("N94").Value = -2385.1666
# Row 134
$ws.Range("H134").Value = 2664
$ws.Range("I134").Value = 2458.4614
$ws.Range("K134").Value = 7375.3842
$ws.Range("M134").Value = -4840.3842

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 120
$ws.Range("I7").Value = 53.333332
$ws.Range("J7").Value = 153.33333
$ws.Range("K7").Value = 53.333332
$ws.Range("L7").Value = 153.33333
$ws.Range("M7").Value = 59.666668
$ws.Range("N7").Value = -379.33333
# Row 31
$ws.Range("H31").Value = 2289.5112
$ws.Range("I31").Value = 2357.1738
$ws.Range("J31").Value = 2218.7727
$ws.Range("K31").Value = 2357.1738
$ws.Range("L31").Value = 2218.7727
$ws.Range("M31").Value = -2062.1738
$ws.Range("N31").Value = -2808.7727
# Row 34
$ws.Range("H34").Value = 2289.5112
$ws.Range("I34").Value = 2357.1738
$ws.Range("J34").Value = 2218.7727
$ws.Range("K34").Value = 2357.1738
$ws.Range("L34").Value = 2218.7727
$ws.Range("M34").Value = -2155.1738
$ws.Range("N34").Value = -2622.7727
# Row 110
$ws.Range("H110").Value = 26500
$ws.Range("J110").Value = 26500
$ws.Range("L110").Value = 26500
$ws.Range("N110").Value = -34680
# Row 132
$ws.Range("H132").Value = 3350
$ws.Range("I132").Value = 3006.08
$ws.Range("J132").Value = 5499.5
$ws.Range("K132").Value = 9018.24
$ws.Range("L132").Value = 16498.5
$ws.Range("M132").Value = -6488.24
$ws.Range("N132").Value = -21558.5
# Row 134
$ws.Range("H134").Value = 3848.5
$ws.Range("I134").Value = 3298.5
$ws.Range("J134").Value = 4673.5
$ws.Range("K134").Value = 9895.5
$ws.Range("L134").Value = 14020.5
$ws.Range("M134").Value = -7360.5
$ws.Range("N134").Value = -19090.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("I68").Value = 567
$ws.Range("J68").Value = 1065.2858
$ws.Range("K68").Value = 1701
$ws.Range("L68").Value = 3195.8574
$ws.Range("M68").Value = -890
$ws.Range("N68").Value = -4817.857400000001
# Row 71
$ws.Range("I71").Value = 567
$ws.Range("J71").Value = 1065.2858
$ws.Range("K71").Value = 5103
$ws.Range("L71").Value = 9587.572200000001
$ws.Range("M71").Value = -1047
$ws.Range("N71").Value = -17699.5722
# Row 92
$ws.Range("H92").Value = 708.6667
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 663
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 1989
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -4485
# Row 120
$ws.Range("H120").Value = 5702
$ws.Range("I120").Value = 5702
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 17106
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -12268
$ws.Range("N120").Value = $null

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 20286.666
$ws.Range("I5").Value = 850
$ws.Range("K5").Value = 850
$ws.Range("M5").Value = -738
# Row 122
$ws.Range("H122").Value = 2697
$ws.Range("I122").Value = 2596
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7788
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5338
$ws.Range("N122").Value = -13900
# Row 132
$ws.Range("H132").Value = 2548.6287
$ws.Range("I132").Value = 2008.16
$ws.Range("J132").Value = 3899.8
$ws.Range("K132").Value = 6024.48
$ws.Range("L132").Value = 11699.4
$ws.Range("M132").Value = -3494.48
$ws.Range("N132").Value = -16759.4

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3455.5715
$ws.Range("I7").Value = 1529.6666
$ws.Range("K7").Value = 1529.6666
$ws.Range("M7").Value = -1417.6666
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = $null
# Row 61
$ws.Range("H61").Value = 19366.5
$ws.Range("I61").Value = 22391.8
$ws.Range("J61").Value = 4240
$ws.Range("K61").Value = 22391.8
$ws.Range("L61").Value = 4240
$ws.Range("M61").Value = -22189.8
$ws.Range("N61").Value = -4644
# Row 113
$ws.Range("H113").Value = 19366.5
$ws.Range("I113").Value = 22391.8
$ws.Range("J113").Value = 4240
$ws.Range("K113").Value = 22391.8
$ws.Range("L113").Value = 4240
$ws.Range("M113").Value = -20221.8
$ws.Range("N113").Value = -8580
# Row 126
$ws.Range("H126").Value = 3455.5715
$ws.Range("I126").Value = 1529.6666
$ws.Range("K126").Value = 4588.9998
$ws.Range("M126").Value = -2118.9998
# Row 132
$ws.Range("H132").Value = 5361.478
$ws.Range("I132").Value = 5546.8823
$ws.Range("J132").Value = 4836.1665
$ws.Range("K132").Value = 16640.6469
$ws.Range("L132").Value = 14508.4995
$ws.Range("M132").Value = -14110.6469
$ws.Range("N132").Value = -19568.4995
# Row 136
$ws.Range("H136").Value = 1354.7858
$ws.Range("I136").Value = 906.0909
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2718.2727
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -168.2727
$ws.Range("N136").Value = -14100

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4689.4287
$ws.Range("I62").Value = 5189.8
$ws.Range("J62").Value = 4411.4443
$ws.Range("K62").Value = 5189.8
$ws.Range("L62").Value = 4411.4443
$ws.Range("M62").Value = -4565.8
$ws.Range("N62").Value = -5659.4443
# Row 65
$ws.Range("H65").Value = 4689.4287
$ws.Range("I65").Value = 5189.8
$ws.Range("J65").Value = 4411.4443
$ws.Range("K65").Value = 25949
$ws.Range("L65").Value = 22057.2215
$ws.Range("M65").Value = -22829
$ws.Range("N65").Value = -28297.2215
# Row 113
$ws.Range("H113").Value = 999.4231
$ws.Range("I113").Value = 741.9091
$ws.Range("J113").Value = 1188.2667
$ws.Range("K113").Value = 2225.7273
$ws.Range("L113").Value = 3564.800099999999
$ws.Range("M113").Value = -55.72730000000001
$ws.Range("N113").Value = -7904.800099999999
# Row 126
$ws.Range("H126").Value = 14332.75
$ws.Range("I126").Value = 17220.889
$ws.Range("J126").Value = 5668.3335
$ws.Range("K126").Value = 51662.667
$ws.Range("L126").Value = 17005.0005
$ws.Range("M126").Value = -49192.667
$ws.Range("N126").Value = -21945.0005
# Row 132
$ws.Range("H132").Value = 2382.6538
$ws.Range("I132").Value = 1863.5122
$ws.Range("J132").Value = 4317.636
$ws.Range("K132").Value = 5590.536599999999
$ws.Range("L132").Value = 12952.908
$ws.Range("M132").Value = -3060.536599999999
$ws.Range("N132").Value = -18012.908
